# Propuesta trabajo final.docx - split/merge runs inside two OMath (equation)
# zones.
#
#  1) In the equation  PIB_mun = f(X) + u  the single run "=f" is split
#     into two separate runs "=" and "f", and the single run "+u" is
#     split into two separate runs "+" and "u".
#  2) In the equation  pi*b_cons = pobl_tot + areaoficialkm2 + ...  the
#     "b" (with an "cons" subscript, via an m:sSub) plus the trailing
#     run-per-token breakdown of the right-hand side are collapsed into
#     a plain run "b_cons" followed by a single run holding the whole
#     "=pobl_tot+areaoficialkm2+discapital+g_cap+finan_credito+
#     vrf_peq_productor+lights_mean" right-hand side.
#
# Word's math runs are not reachable through Range.Text based Find &
# Replace (OMath glyphs get auto-italicised to separate Unicode math
# alphanumeric code points that don't correspond to the underlying plain
# ASCII stored in <m:t>), so the OMath objects are addressed directly
# through Document.OMaths and their content is swapped out wholesale with
# Range.InsertXML, which replaces the contents of the exact range it is
# called on.

$d = $word.ActiveDocument

# --- Equation 2: "PIB_mun=f(X)+u" -> split "=f" and "+u" into singleton runs
$om2 = $d.OMaths.Item(2)
$r2 = $om2.Range.Duplicate

$xml2 = '<m:oMathPara>' +
  '<m:oMath>' +
    '<m:sSub>' +
      '<m:sSubPr>' +
        '<m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr>' +
      '</m:sSubPr>' +
      '<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>PIB</m:t></m:r></m:e>' +
      '<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>mun</m:t></m:r></m:sub>' +
    '</m:sSub>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>f</m:t></m:r>' +
    '<m:d>' +
      '<m:dPr>' +
        '<m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr>' +
      '</m:dPr>' +
      '<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>X</m:t></m:r></m:e>' +
    '</m:d>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>+</m:t></m:r>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>u</m:t></m:r>' +
  '</m:oMath>' +
'</m:oMathPara>'

$r2.InsertXML($xml2)

# --- Equation 3: drop the b_cons m:sSub + per-token runs -> "b_cons" run
#     and one merged right-hand-side run
$om3 = $d.OMaths.Item(3)
$r3 = $om3.Range.Duplicate

$xml3 = '<m:oMathPara>' +
  '<m:oMath>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>pi</m:t></m:r>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>b_cons</m:t></m:r>' +
    '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>=pobl_tot+areaoficialkm2+discapital+g_cap+finan_credito+vrf_peq_productor+lights_mean</m:t></m:r>' +
  '</m:oMath>' +
'</m:oMathPara>'

$r3.InsertXML($xml3)

Write-Output "OK"
